$wb = $excel.ActiveWorkbook

# --- 1. Update selection on the "survey" sheet (A10:XFD10 -> E14) ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("E14").Select()

# --- 2. Add the new "properties" sheet after the last existing sheet ("choices") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$props = $wb.Worksheets.Add($null, $lastSheet)
$props.Name = "properties"

# Header row
$props.Range("A1").Value = "partition"
$props.Range("B1").Value = "aspect"
$props.Range("C1").Value = "key"
$props.Range("D1").Value = "type"
$props.Range("E1").Value = "value"

# Data row (written in this column order to mirror the original shared-string layout)
$props.Range("C2").Value = "colOrder"
$props.Range("B2").Value = "default"
$props.Range("A2").Value = "Table"
$props.Range("E2").Value = '["Date_and_Time","plot_name","soil_condition","Sample"]'
$props.Range("D2").Value = "array"

# Selection on the new sheet, which also becomes the active tab/sheet
$props.Range("D3").Select()
